$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header cells for new columns AD, AE, AF
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Copy style from an existing header cell (e.g. AC1) to the new headers
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Fill in Wins / Losses / Ties values for each data row (2-46)
for ($row = 2; $row -le 46; $row++) {
    $ws.Cells.Item($row, 30).Value = 66   # AD = column 30
    $ws.Cells.Item($row, 31).Value = 96   # AE = column 31
    $ws.Cells.Item($row, 32).Value = 0    # AF = column 32
}
